$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The coefficient table (B2:J9) was displaying raw full-precision numbers;
# round the display to two decimal places (built-in format "0.00", numFmtId 2)
# while leaving the header row/label column formatting untouched.
$ws.Range("B2:J9").NumberFormat = "0.00"

# Columns were resized to fit their (now 2-decimal) contents. Group the
# columns the same way Excel's own best-fit pass grouped them: B:E share one
# width, F:G share a slightly wider one, and H:J (the largest numbers) share
# the widest.
$ws.Range("B1:E1").EntireColumn.ColumnWidth = 8.1666666666667
$ws.Range("F1:G1").EntireColumn.ColumnWidth = 8.6666666666667
$ws.Range("H1:J1").EntireColumn.ColumnWidth = 10.3333333333333

# The active selection moved from J3 to L3.
$ws.Range("L3").Select()
